$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Avinash"
$ws.Range("B1").Value = ".Net Intern"
$ws.Range("C1").Value = "First line added"

$ws.Range("D1").Select()
